$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.818.71"
$ws.Range("E2").Value = "  -3.06%  "

$ws.Range("D3").Value = "2.615.51"
$ws.Range("E3").Value = "  -1.72%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'575.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.01%  "

$ws.Range("D6").Value = "'156.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.09%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").Value = "2.613.14"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("E10").Value = "  -6.63%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").Value = "'0.381"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.95%  "

$ws.Range("E13").Value = "  -0.31%  "

$ws.Range("D14").Value = "'28.27"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.94%  "

$ws.Range("D15").Value = "3.084.34"
$ws.Range("E15").Value = "  -1.76%  "

$ws.Range("D16").Value = "'0.0000181"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -7.75%  "

$ws.Range("D17").Value = "63.655.25"
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").Value = "2.618.09"
$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").Value = "'12.04"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.40%  "

$ws.Range("D20").Value = "'7.60"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "'4.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.51%  "

$ws.Range("D22").Value = "'344.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.75%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").Value = "'67.54"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.59%  "

$ws.Range("D25").Value = "'1.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("E26").Value = "  -3.41%  "

$ws.Range("D27").Value = "'597.37"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").Value = "'9.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.96%  "

$ws.Range("D29").Value = "'1.58"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.77%  "

$ws.Range("D30").Value = "'0.162"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").Value = "'0.997"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.27%  "

$ws.Range("D32").Value = "'7.92"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.86%  "

$ws.Range("E33").Value = "  -4.21%  "

$ws.Range("E34").Value = "  -3.97%  "

$ws.Range("D35").Value = "'6.61"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.95%  "

$ws.Range("D36").Value = "'5.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.50%  "

$ws.Range("D37").Value = "'0.403"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.62%  "

$ws.Range("D38").Value = "'19.74"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.16%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'154.13"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'1.88"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.38%  "

# Rows 42 and 43 swap coin data (dogwifhat moves above USDe) with updated values
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "'41.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.39%  "

$ws.Range("D45").Value = "'157.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").Value = "'23.86"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").Value = "'3.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.88%  "

$ws.Range("D48").Value = "'0.0591"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.20%  "

$ws.Range("E49").Value = "  -1.92%  "

$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("E51").Value = "  -5.07%  "

